{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst para = paragraphs.items[0];\npara.load(\"text\");\nawait context.sync();\n\n// Replace the whole paragraph's text (currently split across several\n// runs spelling \"Buenos d\u00edas\") with the new single run of text.\npara.insertText(\"Versi\u00f3n con errores\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The paragraph currently reads \"Buenos d\u00edas\" (split across several\n# runs). Replace the whole paragraph text with the new text while\n# keeping the paragraph mark (and its formatting) intact.\n$para = $d.Paragraphs(1)\n$r = $para.Range\n$r.MoveEnd(1, -1) | Out-Null\n$r.Text = \"Versi\u00f3n con errores\"\n"}
